$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# input heuristics (single-dot decimals). Force Text format first so the
# stored value stays a string, matching the source data's inline-string cells.
$textForced = @("D5", "D6", "D8", "D10", "D12", "D13", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textForced) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values (order follows the sheet, top to bottom).
$ws.Range("D2").Value = "59.970.72"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.584.43"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "563.28"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "141.38"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "2.598.72"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "0.367"
$ws.Range("E12").Value = "  +6.25%  "
$ws.Range("D13").Value = "0.151"
$ws.Range("E13").Value = "  -5.60%  "
$ws.Range("D14").Value = "3.041.17"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "59.978.76"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "23.21"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "2.592.53"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +8.39%  "
$ws.Range("D20").Value = "4.62"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").Value = "344.44"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "6.94"
$ws.Range("E22").Value = "  +8.42%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "0.529"
$ws.Range("E24").Value = "  +15.84%  "
$ws.Range("D25").Value = "62.70"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "7.59"
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +6.32%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").Value = "161.03"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "19.36"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("E36").Value = "  +8.48%  "
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").Value = "37.54"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "0.856"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").Value = "3.79"
$ws.Range("D42").Value = "292.37"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "137.58"
$ws.Range("E43").Value = "  +4.64%  "
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0975"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.601"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "0.0542"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "19.36"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "10.66"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.83"
$ws.Range("E51").Value = "  +6.91%  "
